$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "A3"  = -21.9413
    "C3"  = -11.8687
    "A21" = -20.17409999999998
    "A23" = -20.29349999999998
    "C24" = -12.77309999999999
    "A25" = -21.76590000000001
    "B27" = 5.901400000000003
    "B31" = 5.331199999999999
    "B39" = 9.804800000000006
    "B48" = 5.442200000000003
    "B51" = 5.455999999999997
    "B52" = 5.348700000000001
    "A53" = -21.89889999999999
    "B55" = 5.874399999999997
    "B56" = 4.929799999999998
    "A57" = -22.03080000000001
    "B57" = 5.162900000000002
    "C57" = -12.35459999999999
    "A59" = -22.2938
    "C61" = -13.60129999999999
    "A69" = -21.58899999999999
    "C70" = -11.5894
    "B73" = 8.524299999999998
    "A79" = -20.54420000000001
    "A83" = -21.92219999999999
    "C86" = -13.0989
    "B89" = 5.129499999999994
    "B90" = 5.702200000000001
    "A93" = -21.55560000000001
    "C98" = -11.8981
    "C100" = -12.39569999999999
    "C102" = -13.3729
}

foreach ($addr in $changes.Keys) {
    $ws.Range($addr).Value = $changes[$addr]
}
